$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) column names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize connector words (de/del/el/la/los/las/y) in state and municipality names
$ws.Range("B7").Value = "Pabellón De Arteaga"
$ws.Range("B8").Value = "Rincón De Romos"
$ws.Range("B9").Value = "San Francisco De Los Romo"
$ws.Range("B10").Value = "San José De Gracia"
$ws.Range("B15").Value = "Playas De Rosarito"
$ws.Range("B33").Value = "Bejucal De Ocampo"
$ws.Range("B56").Value = "Ocozocoautla De Espinosa"
$ws.Range("B62").Value = "Salto De Agua"
$ws.Range("B98").Value = "Guadalupe Y Calvo"
$ws.Range("B101").Value = "Hidalgo Del Parral"
$ws.Range("B120").Value = "San Francisco De Borja"
$ws.Range("B121").Value = "San Francisco Del Oro"
$ws.Range("B125").Value = "Valle De Zaragoza"
$ws.Range("B154").Value = "Villa De Álvarez"
$ws.Range("A156").Value = "Ciudad De México"
$ws.Range("B160").Value = "Cuajimalpa De Morelos"
$ws.Range("B175").Value = "Coneto De Comonfort"
$ws.Range("B189").Value = "Nombre De Dios"
$ws.Range("B199").Value = "San Juan Del Río"
$ws.Range("A207").Value = "Estado De México"
$ws.Range("B207").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B210").Value = "Almoloya De Juárez"
$ws.Range("B215").Value = "Atizapán De Zaragoza"
$ws.Range("B220").Value = "Chapa De Mota"
$ws.Range("B223").Value = "Coacalco De Berriozábal"
$ws.Range("B229").Value = "Ecatepec De Morelos"
$ws.Range("B233").Value = "Ixtapan De La Sal"
$ws.Range("B245").Value = "Naucalpan De Juárez"
$ws.Range("B256").Value = "San Felipe Del Progreso"
$ws.Range("B257").Value = "San Martín De Las Pirámides"
$ws.Range("B267").Value = "Tenango Del Aire"
$ws.Range("B268").Value = "Tenango Del Valle"
$ws.Range("B275").Value = "Tlalnepantla De Baz"
$ws.Range("B279").Value = "Valle De Chalco Solidaridad"
$ws.Range("B280").Value = "Villa De Allende"
$ws.Range("B291").Value = "San Miguel De Allende"
$ws.Range("B292").Value = "Apaseo El Alto"
$ws.Range("B293").Value = "Apaseo El Grande"
$ws.Range("B299").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B303").Value = "Jaral Del Progreso"
$ws.Range("B311").Value = "Purísima Del Rincón"
$ws.Range("B316").Value = "San Francisco Del Rincón"
$ws.Range("B318").Value = "San Luis De La Paz"
$ws.Range("B319").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B321").Value = "Silao De La Victoria"
$ws.Range("B325").Value = "Valle De Santiago"
$ws.Range("B331").Value = "Acapulco De Juárez"
$ws.Range("B333").Value = "Ajuchitlán Del Progreso"
$ws.Range("B334").Value = "Alcozauca De Guerrero"
$ws.Range("B338").Value = "Atoyac De Álvarez"
$ws.Range("B339").Value = "Ayutla De Los Libres"
$ws.Range("B342").Value = "Buenavista De Cuéllar"
$ws.Range("B343").Value = "Chilapa De Álvarez"
$ws.Range("B344").Value = "Chilpancingo De Los Bravo"
$ws.Range("B345").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B349").Value = "Coyuca De Benítez"
$ws.Range("B350").Value = "Coyuca De Catalán"
$ws.Range("B354").Value = "Cuetzala Del Progreso"
$ws.Range("B355").Value = "Cutzamala De Pinzón"
$ws.Range("B360").Value = "Huitzuco De Los Figueroa"
$ws.Range("B361").Value = "Iguala De La Independencia"
$ws.Range("B363").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B364").Value = "Zihuatanejo De Azueta"
$ws.Range("B366").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B380").Value = "Taxco De Alarcón"
$ws.Range("B382").Value = "Técpan De Galeana"
$ws.Range("B384").Value = "Tepecoacuilco De Trujano"
$ws.Range("B386").Value = "Tixtla De Guerrero"
$ws.Range("B389").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B390").Value = "Tlapa De Comonfort"
$ws.Range("B403").Value = "Atotonilco De Tula"
$ws.Range("B404").Value = "Atotonilco El Grande"
$ws.Range("B409").Value = "Cuautepec De Hinojosa"
$ws.Range("B412").Value = "Huasca De Ocampo"
$ws.Range("B415").Value = "Huejutla De Reyes"
$ws.Range("B418").Value = "Jacala De Ledezma"
$ws.Range("B421").Value = "Mineral Del Monte"
$ws.Range("B422").Value = "Mixquiahuala De Juárez"
$ws.Range("B423").Value = "Molango De Escamilla"
$ws.Range("B425").Value = "Nopala De Villagrán"
$ws.Range("B426").Value = "Omitlán De Juárez"
$ws.Range("B427").Value = "Pachuca De Soto"
$ws.Range("B428").Value = "Progreso De Obregón"
$ws.Range("B434").Value = "Santiago De Anaya"
$ws.Range("B435").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B438").Value = "Tenango De Doria"
$ws.Range("B442").Value = "Tezontepec De Aldama"
$ws.Range("B447").Value = "Tula De Allende"
$ws.Range("B448").Value = "Tulancingo De Bravo"
$ws.Range("B449").Value = "Villa De Tezontepec"
$ws.Range("B452").Value = "Zacualtipán De Ángeles"
$ws.Range("B453").Value = "Zapotlán De Juárez"
$ws.Range("B458").Value = "Acatlán De Juárez"
$ws.Range("B459").Value = "Ahualulco De Mercado"
$ws.Range("B464").Value = "Atemajac De Brizuela"
$ws.Range("B467").Value = "Atotonilco El Alto"
$ws.Range("B469").Value = "Autlán De Navarro"
$ws.Range("B481").Value = "Concepción De Buenos Aires"
$ws.Range("B482").Value = "Cuautitlán De García Barragán"
$ws.Range("B491").Value = "Encarnación De Díaz"
$ws.Range("B498").Value = "Huejuquilla El Alto"
$ws.Range("B499").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B500").Value = "Ixtlahuacán Del Río"
$ws.Range("B504").Value = "Jilotlán De Los Dolores"
$ws.Range("B510").Value = "La Manzanilla De La Paz"
$ws.Range("B511").Value = "Lagos De Moreno"
$ws.Range("B518").Value = "Ojuelos De Jalisco"
$ws.Range("B523").Value = "San Cristóbal De La Barranca"
$ws.Range("B524").Value = "San Diego De Alejandría"
$ws.Range("B526").Value = "San Juan De Los Lagos"
$ws.Range("B527").Value = "San Juanito De Escobedo"
$ws.Range("B529").Value = "San Martín De Bolaños"
$ws.Range("B531").Value = "San Miguel El Alto"
$ws.Range("B532").Value = "San Sebastián Del Oeste"
$ws.Range("B533").Value = "Santa María De Los Ángeles"
$ws.Range("B534").Value = "Santa María Del Oro"
$ws.Range("B537").Value = "Talpa De Allende"
$ws.Range("B538").Value = "Tamazula De Gordiano"
$ws.Range("B541").Value = "Techaluta De Montenegro"
$ws.Range("B545").Value = "Teocuitatlán De Corona"
$ws.Range("B546").Value = "Tepatitlán De Morelos"
$ws.Range("B549").Value = "Tizapán El Alto"
$ws.Range("B550").Value = "Tlajomulco De Zúñiga"
$ws.Range("B562").Value = "Unión De San Antonio"
$ws.Range("B563").Value = "Unión De Tula"
$ws.Range("B567").Value = "Yahualica De González Gallo"
$ws.Range("B568").Value = "Zacoalco De Torres"
$ws.Range("B571").Value = "Zapotitlán De Vadillo"
$ws.Range("B572").Value = "Zapotlán Del Rey"
$ws.Range("B573").Value = "Zapotlán El Grande"
$ws.Range("B596").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B598").Value = "Cojumatlán De Régules"
$ws.Range("B661").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B693").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B696").Value = "Puente De Ixtla"
$ws.Range("B702").Value = "Tlaltizapán De Zapata"
$ws.Range("B708").Value = "Zacualpan De Amilpas"
$ws.Range("B712").Value = "Amatlán De Cañas"
$ws.Range("B713").Value = "Bahía De Banderas"
$ws.Range("B717").Value = "Ixtlán Del Río"
$ws.Range("B724").Value = "Santa María Del Oro"
$ws.Range("B733").Value = "Lampazos De Naranjo"
$ws.Range("B736").Value = "San Nicolás De Los Garza"
$ws.Range("B739").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B743").Value = "Ayoquezco De Aldama"
$ws.Range("B744").Value = "Chalcatongo De Hidalgo"
$ws.Range("B745").Value = "Coicoyán De Las Flores"
$ws.Range("B748").Value = "Fresnillo De Trujano"
$ws.Range("B750").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B751").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B752").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B753").Value = "Ixtlán De Juárez"
$ws.Range("B754").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B761").Value = "Mariscala De Juárez"
$ws.Range("B763").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B764").Value = "Mixistlán De La Reforma"
$ws.Range("B766").Value = "Nejapa De Madero"
$ws.Range("B767").Value = "Oaxaca De Juárez"
$ws.Range("B768").Value = "Ocotlán De Morelos"
$ws.Range("B769").Value = "Pinotepa De Don Luis"
$ws.Range("B771").Value = "Putla Villa De Guerrero"
$ws.Range("B778").Value = "San Antonio De La Cal"
$ws.Range("B792").Value = "San José Del Progreso"
$ws.Range("B795").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B800").Value = "San Juan Del Estado"
$ws.Range("B835").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B845").Value = "Santa Cruz Tacache De Mina"
$ws.Range("B851").Value = "Santa María Del Tule"
$ws.Range("B885").Value = "Tataltepec De Valdés"
$ws.Range("B886").Value = "Teococuilco De Marcos Pérez"
$ws.Range("B887").Value = "Teotitlán De Flores Magón"
$ws.Range("B888").Value = "Teotitlán Del Valle"
$ws.Range("B889").Value = "Tepelmeme Villa De Morelos"
$ws.Range("B890").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B891").Value = "Tlacolula De Matamoros"
$ws.Range("B893").Value = "Villa De Chilapa De Díaz"
$ws.Range("B894").Value = "Villa De Etla"
$ws.Range("B895").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B896").Value = "Villa De Tututepec"
$ws.Range("B897").Value = "Villa De Zaachila"
$ws.Range("B899").Value = "Villa Sola De Vega"
$ws.Range("B900").Value = "Zapotitlán Del Río"
$ws.Range("B901").Value = "Zimatlán De Álvarez"
$ws.Range("B922").Value = "Chalchicomula De Sesma"
$ws.Range("B931").Value = "Chila De La Sal"
$ws.Range("B937").Value = "Cuapiaxtla De Madero"
$ws.Range("B940").Value = "Cuayuca De Andrade"
$ws.Range("B950").Value = "Huehuetlán El Chico"
$ws.Range("B951").Value = "Huehuetlán El Grande"
$ws.Range("B957").Value = "Ixcamilpa De Guerrero"
$ws.Range("B959").Value = "Izúcar De Matamoros"
$ws.Range("B966").Value = "Los Reyes De Juárez"
$ws.Range("B967").Value = "Mazapiltepec De Juárez"
$ws.Range("B975").Value = "Palmar De Bravo"
$ws.Range("B992").Value = "San Nicolás De Los Ranchos"
$ws.Range("B995").Value = "San Salvador El Seco"
$ws.Range("B996").Value = "San Salvador El Verde"
$ws.Range("B1001").Value = "Tecali De Herrera"
$ws.Range("B1009").Value = "Tepanco De López"
$ws.Range("B1010").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B1014").Value = "Tepexi De Rodríguez"
$ws.Range("B1016").Value = "Tepeyahualco De Cuauhtémoc"
$ws.Range("B1017").Value = "Tetela De Ocampo"
$ws.Range("B1022").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B1031").Value = "Tuzamapan De Galeana"
$ws.Range("B1035").Value = "Xayacatlán De Bravo"
$ws.Range("B1049").Value = "Amealco De Bonfil"
$ws.Range("B1050").Value = "Cadereyta De Montes"
$ws.Range("B1054").Value = "Jalpan De Serra"
$ws.Range("B1057").Value = "Pinal De Amoles"
$ws.Range("B1060").Value = "San Juan Del Río"
$ws.Range("B1070").Value = "Cerro De San Pedro"
$ws.Range("B1072").Value = "Ciudad Del Maíz"
$ws.Range("B1079").Value = "Mexquitic De Carmona"
$ws.Range("B1084").Value = "Santa María Del Río"
$ws.Range("B1090").Value = "Villa De Arriaga"
$ws.Range("B1091").Value = "Villa De Ramos"
$ws.Range("B1092").Value = "Villa De Reyes"
$ws.Range("B1136").Value = "Nacozari De García"
$ws.Range("B1158").Value = "Jalpa De Méndez"
$ws.Range("B1179").Value = "Soto La Marina"
$ws.Range("B1184").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B1188").Value = "Contla De Juan Cuamatzi"
$ws.Range("B1193").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B1195").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B1198").Value = "San Pablo Del Monte"
$ws.Range("B1199").Value = "Tepetitla De Lardizábal"
$ws.Range("B1200").Value = "Tetla De La Solidaridad"
$ws.Range("B1214").Value = "Amatlán De Los Reyes"
$ws.Range("B1220").Value = "Boca Del Río"
$ws.Range("B1235").Value = "Cosamaloapan De Carpio"
$ws.Range("B1236").Value = "Cosautlán De Carvajal"
$ws.Range("B1247").Value = "Hueyapan De Ocampo"
$ws.Range("B1248").Value = "Huiloapan De Cuauhtémoc"
$ws.Range("B1251").Value = "Ixhuatlán De Madero"
$ws.Range("B1252").Value = "Ixhuatlán Del Sureste"
$ws.Range("B1262").Value = "Lerdo De Tejada"
$ws.Range("B1266").Value = "Martínez De La Torre"
$ws.Range("B1268").Value = "Medellín De Bravo"
$ws.Range("B1278").Value = "Paso De Ovejas"
$ws.Range("B1279").Value = "Paso Del Macho"
$ws.Range("B1283").Value = "Poza Rica De Hidalgo"
$ws.Range("B1291").Value = "Sayula De Alemán"
$ws.Range("B1294").Value = "Soledad De Doblado"
$ws.Range("B1342").Value = "Concepción Del Oro"
$ws.Range("B1344").Value = "El Plateado De Joaquín Amaro"
$ws.Range("B1358").Value = "Mezquital Del Oro"
$ws.Range("B1362").Value = "Moyahua De Estrada"
$ws.Range("B1363").Value = "Nochistlán De Mejía"
$ws.Range("B1364").Value = "Noria De Ángeles"
$ws.Range("B1374").Value = "Teúl De González Ortega"
$ws.Range("B1375").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B1377").Value = "Villa De Cos"

# Refresh percentage values that shift by a tiny floating point rounding amount
$ws.Range("D109").Value2 = 0.0009566326530612244
$ws.Range("D191").Value2 = 0.0009566326530612244
$ws.Range("D205").Value2 = 0.0009566326530612244
$ws.Range("D240").Value2 = 0.0009566326530612244
$ws.Range("D349").Value2 = 0.0009566326530612244
$ws.Range("D478").Value2 = 0.0009566326530612244
$ws.Range("D487").Value2 = 0.0009566326530612244
$ws.Range("D636").Value2 = 0.0009566326530612244
$ws.Range("D705").Value2 = 0.0009566326530612244
$ws.Range("D710").Value2 = 0.0009566326530612244
$ws.Range("D750").Value2 = 0.0009566326530612244
$ws.Range("D927").Value2 = 0.0009566326530612244
$ws.Range("D1114").Value2 = 0.0009566326530612244
$ws.Range("D1124").Value2 = 0.0009566326530612244
$ws.Range("D1265").Value2 = 0.0009566326530612244
$ws.Range("D1314").Value2 = 0.0009566326530612244
$ws.Range("D1336").Value2 = 0.0009566326530612244

# Remove footer/metadata rows 1384-1388 (source notes no longer needed)
$ws.Rows("1384:1388").Delete()

